$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change D1 header from "Order" to "Step"
$ws.Range("D1").Value = "Step"

# Update the selected cell to D2, matching the author's new selection
$ws.Range("D2").Select()
